$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data
$ws.Range("D2").Value = '93.503.66'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '3.426.76'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''234.14'
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").Value = '''621.76'
$ws.Range("E6").Value = '  -2.53%  '
$ws.Range("D7").Value = '''1.40'
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("D8").Value = '''0.397'
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").Value = '3.425.83'
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '''43.07'
$ws.Range("E12").Value = '  +4.61%  '
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '''6.29'
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = '93.277.48'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '4.075.18'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '''8.22'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '3.432.18'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '''18.10'
$ws.Range("E20").Value = '  +4.56%  '
$ws.Range("D21").Value = '''11.70'
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = '''3.40'
$ws.Range("E22").Value = '  +5.43%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '''502.91'
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("E25").Value = '  +3.22%  '
$ws.Range("E26").Value = '  -2.67%  '
$ws.Range("D27").Value = '''95.09'
$ws.Range("E27").Value = '  +4.77%  '
$ws.Range("D28").Value = '''11.99'
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("D29").Value = '3.611.22'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").Value = '''11.47'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  +2.38%  '
$ws.Range("D33").Value = '''2.76'
$ws.Range("E33").Value = '  +1.86%  '
$ws.Range("D34").Value = '''0.991'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("D37").Value = '''0.553'
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("D38").Value = '''558.45'
$ws.Range("E38").Value = '  +4.30%  '
$ws.Range("D39").Value = '''7.51'
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").Value = '''0.913'
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").Value = '''23.70'
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0411'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = '''5.52'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '''53.66'
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("D50").Value = '''2.14'
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = '''8.12'
$ws.Range("E51").Value = '  +1.74%  '

# Reset number format/style on cells forced to text via leading apostrophe
# so no stray style index remains on these cells
$resetCells = @("D5","D6","D7","D8","D12","D14","D18","D20","D21","D22","D23","D27","D28","D30","D33","D34","D37","D38","D39","D43","D45","D47","D48","D49","D50","D51")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}

